# Normalize the ordering of names/emails in the "Recorded By" column (G).
# The values are comma-separated lists of the same recorder identities;
# this re-orders them (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, backup@backdoor.com, system" = "System, system, backup@backdoor.com"
    "dnasr281@gmail.com, System" = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text
    $replacement = $map[$current]
    if ($replacement) {
        $cell.Value = $replacement
    }
}
